$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: this COM host's "case-sensitive" comparison operators (-ceq/-cne/
# -clike/-cmatch) do not actually distinguish case (observed empirically),
# so case-sensitive exact-string matching is done manually via character
# codes instead.
function Test-ExactMatch($a, $b) {
    if ($a.Length -ne $b.Length) { return $false }
    $chars_a = $a.ToCharArray()
    $chars_b = $b.ToCharArray()
    for ($i = 0; $i -lt $chars_a.Length; $i++) {
        if ([int]$chars_a[$i] -ne [int]$chars_b[$i]) {
            return $false
        }
    }
    return $true
}

# Determine the extent of the used range (data starts at row 2; row 1 is the header).
$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = @($val -split ",\s*")

        $hasSystem = $false
        foreach ($p in $parts) {
            if (Test-ExactMatch $p "System") { $hasSystem = $true }
        }

        $firstIsSystem = Test-ExactMatch $parts[0] "System"

        # Only reorder when "System" (exact case) appears among the comma
        # separated names and isn't already the first entry. A distinct,
        # lowercase "system" entry (if present) must be left alone/in place.
        if ($hasSystem -and -not $firstIsSystem) {
            $rest = @()
            foreach ($p in $parts) {
                if (-not (Test-ExactMatch $p "System")) {
                    $rest += $p
                }
            }
            $newParts = @("System") + $rest
            $newVal = [string]::Join(", ", $newParts)
            $cell.Value2 = $newVal
        }
    }
}
